# IST price update 2025-12-23 04:00
# A new price-check timestamp column is inserted as the new column B (most
# recent check), shifting every existing data column one slot to the right
# (old B..AE -> new C..AF). The freshly inserted column is seeded with the
# new check's timestamp header and carries forward the latest known price
# for each SKU (i.e. the same price that is now sitting in column C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Figure out how many data rows currently exist (row 1 = header row).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

# Insert a new column before column B; this shifts old B:AE -> C:AF and
# widens the used range to A1:AF26.
$ws.Columns("B").Insert()

# Match the column width used by every other price column (raw OOXML
# width of 21 characters).
$ws.Range("B1").ColumnWidth = 20.17

# New timestamp header for the freshly inserted (most recent) check.
$ws.Range("B1").Value = "2025-12-23 09:25"
$ws.Range("B1").Style = $ws.Range("C1").Style

# Carry the latest known price (now duplicated in column C after the
# shift) forward into the brand-new column B for every SKU row.
for ($r = 2; $r -le $lastRow; $r++) {
    $src = $ws.Cells.Item($r, 3)
    $dst = $ws.Cells.Item($r, 2)
    if ($src.Value -eq $null) {
        $dst.Value = ""
    } else {
        $dst.Value = $src.Value
    }
}
